## Add 15 new session rows (S9 / Alejandro, 2013-07-02) to the watch-ERP
## dataset sheet, following the same layout as the existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# --- 1. Prepare rows 122-136 with the same formatting (styles) as row 121 ---
$ws.Range("A121:H121").Copy()
$ws.Range("A122:H136").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- 2. Fill in the new data -------------------------------------------------
# Columns: A subjectTag, B subjectName, C date, D sessionDirectory,
#          E fileName, F condition, G frequency, H run

$subjectTag = "S9"
$subjectName = "Alejandro"
$sessionDate = 41457
$sessionDirectory = "2013-07-02-alejandro"

$rows = @(
    @("2013-07-02-10-41-52", "oddball",       0,    1),
    @("2013-07-02-10-51-04", "hybrid-12Hz",   12,   1),
    @("2013-07-02-10-58-51", "hybrid-12Hz",   12,   2),
    @("2013-07-02-11-05-55", "oddball",       0,    2),
    @("2013-07-02-11-14-05", "hybrid-15Hz",   15,   1),
    @("2013-07-02-11-37-37", "hybrid-8-57Hz", 8.57, 1),
    @("2013-07-02-11-45-28", "hybrid-10Hz",   10,   1),
    @("2013-07-02-11-52-29", "oddball",       0,    3),
    @("2013-07-02-11-59-14", "hybrid-15Hz",   15,   2),
    @("2013-07-02-12-06-08", "hybrid-10Hz",   10,   2),
    @("2013-07-02-12-18-34", "hybrid-8-57Hz", 8.57, 2),
    @("2013-07-02-12-24-49", "hybrid-10Hz",   10,   3),
    @("2013-07-02-12-31-13", "hybrid-8-57Hz", 8.57, 3),
    @("2013-07-02-12-37-25", "hybrid-15Hz",   15,   3),
    @("2013-07-02-12-43-59", "hybrid-12Hz",   12,   3)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = 122 + $i
    $item = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $subjectTag
    $ws.Cells.Item($r, 2).Value = $subjectName
    $ws.Cells.Item($r, 3).Value = $sessionDate
    $ws.Cells.Item($r, 4).Value = $sessionDirectory
    $ws.Cells.Item($r, 5).Value = $item[0]
    $ws.Cells.Item($r, 6).Value = $item[1]
    $ws.Cells.Item($r, 7).Value = $item[2]
    $ws.Cells.Item($r, 8).Value = $item[3]
}

# --- 3. Update the view / selection to match the post-edit state -----------
[void]$ws.Range("A124").Select()
